# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to the freshly-scraped counts (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 583
$wsExpo.Range("F4").Value = 1255
$wsExpo.Range("F5").Value = 1087
$wsExpo.Range("F6").Value = 14064
$wsExpo.Range("F7").Value = 15516
$wsExpo.Range("F8").Value = 7
$wsExpo.Range("F9").Value = 46
$wsExpo.Range("F10").Value = 46
$wsExpo.Range("F11").Value = 185
$wsExpo.Range("F12").Value = 24
$wsExpo.Range("F13").Value = 49
$wsExpo.Range("F19").Value = 26
$wsExpo.Range("F20").Value = 1206
$wsExpo.Range("F21").Value = 129
$wsExpo.Range("F22").Value = 64
$wsExpo.Range("F23").Value = 6048
$wsExpo.Range("F24").Value = 958
$wsExpo.Range("F25").Value = 1081
$wsExpo.Range("F26").Value = 5534
$wsExpo.Range("F29").Value = 107
$wsExpo.Range("F30").Value = 2992

# --- Sheet "全部类型" (all types, union of the other sheets) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 583
$wsAll.Range("F5").Value = 1255
$wsAll.Range("F6").Value = 1087
$wsAll.Range("F7").Value = 14064
$wsAll.Range("F8").Value = 15516
$wsAll.Range("F9").Value = 7
$wsAll.Range("F10").Value = 46
$wsAll.Range("F11").Value = 46
$wsAll.Range("F12").Value = 185
$wsAll.Range("F13").Value = 24
$wsAll.Range("F14").Value = 49
$wsAll.Range("F20").Value = 26
$wsAll.Range("F21").Value = 1206
$wsAll.Range("F22").Value = 129
$wsAll.Range("F23").Value = 64
$wsAll.Range("F25").Value = 6048
$wsAll.Range("F26").Value = 958
$wsAll.Range("F27").Value = 1081
$wsAll.Range("F28").Value = 5534
$wsAll.Range("F31").Value = 107
$wsAll.Range("F32").Value = 2998
